$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.033.34"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").Value = "2.549.65"

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "585.04"
$ws.Range("E5").Value = "  +2.34%  "

# Row 6
$ws.Range("D6").Value = "147.10"
$ws.Range("E6").Value = "  -2.48%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -1.09%  "

# Row 9
$ws.Range("E9").Value = "  -0.55%  "

# Row 10
$ws.Range("E10").Value = "  -3.36%  "

# Row 11
$ws.Range("E11").Value = "  -0.19%  "

# Row 12
$ws.Range("E12").Value = "  -1.29%  "

# Row 13
$ws.Range("D13").Value = "27.39"
$ws.Range("E13").Value = "  -3.97%  "

# Row 14
$ws.Range("D14").Value = "3.005.45"
$ws.Range("E14").Value = "  +0.12%  "

# Row 15
$ws.Range("D15").Value = "62.926.19"
$ws.Range("E15").Value = "  -0.54%  "

# Row 16
$ws.Range("E16").Value = "  -0.78%  "

# Row 17
$ws.Range("D17").Value = "2.551.42"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").Value = "11.35"
$ws.Range("E18").Value = "  -2.83%  "

# Row 19
$ws.Range("D19").Value = "336.28"
$ws.Range("E19").Value = "  -1.35%  "

# Row 20
$ws.Range("D20").Value = "4.33"
$ws.Range("E20").Value = "  -1.25%  "

# Row 21
$ws.Range("D21").Value = "6.77"
$ws.Range("E21").Value = "  -1.73%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("E24").Value = "  -0.60%  "

# Row 25
$ws.Range("E25").Value = "  +0.73%  "

# Row 26
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").Value = "1.49"
$ws.Range("E26").Value = "  +0.87%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("D28").Value = "8.37"
$ws.Range("E28").Value = "  -1.26%  "

# Row 29
$ws.Range("D29").Value = "7.49"
$ws.Range("E29").Value = "  +5.84%  "

# Row 30
$ws.Range("D30").Value = "1.93"
$ws.Range("E30").Value = "  +3.47%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0812"
$ws.Range("E31").Value = "  -2.94%  "

# Row 32
$ws.Range("D32").Value = "177.87"
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
$ws.Range("E33").Value = "  -1.06%  "

# Row 34
$ws.Range("D34").Value = "414.90"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35
$ws.Range("D35").Value = "19.14"
$ws.Range("E35").Value = "  -0.13%  "

# Row 36
$ws.Range("E36").Value = "  -1.85%  "

# Row 38
$ws.Range("D38").Value = "4.34"
$ws.Range("E38").Value = "  -2.82%  "

# Row 39
$ws.Range("D39").Value = "1.74"

# Row 40
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("D41").Value = "39.66"
$ws.Range("E41").Value = "  -0.91%  "

# Row 42
$ws.Range("D42").Value = "150.91"
$ws.Range("E42").Value = "  -3.06%  "

# Row 44
$ws.Range("D44").Value = "20.84"
$ws.Range("E44").Value = "  -1.23%  "

# Row 45
$ws.Range("E45").Value = "  +1.31%  "

# Row 46
$ws.Range("E46").Value = "  -1.30%  "

# Row 47
$ws.Range("E47").Value = "  +0.42%  "

# Row 48
$ws.Range("E48").Value = "  -0.39%  "

# Row 49
$ws.Range("D49").Value = "18.27"
$ws.Range("E49").Value = "  -1.94%  "

# Row 50
$ws.Range("E50").Value = "  -6.25%  "

# Row 51
$ws.Range("D51").Value = "11.30"
$ws.Range("E51").Value = "  -0.26%  "
